# Apply the gh-pages data-refresh update to 杭州-漫展信息.xlsx
# Sheets: 1=展览 (Exhibition), 2=演出 (Performance), 3=本地生活 (Local life), 4=全部类型 (All types)

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("展览")
$ws2 = $wb.Worksheets.Item("演出")
$ws3 = $wb.Worksheets.Item("本地生活")
$ws4 = $wb.Worksheets.Item("全部类型")

# ---------------------------------------------------------------------------
# Sheet "展览" (展览) - update "想去人数" (column F) for a number of rows
# ---------------------------------------------------------------------------
$ws1.Range("F2").Value = 1500
$ws1.Range("F3").Value = 837
$ws1.Range("F4").Value = 432
$ws1.Range("F5").Value = 884
$ws1.Range("F7").Value = 7401
$ws1.Range("F10").Value = 1930
$ws1.Range("F11").Value = 5385
$ws1.Range("F12").Value = 559
$ws1.Range("F15").Value = 7450
$ws1.Range("F16").Value = 8761
$ws1.Range("F18").Value = 1133
$ws1.Range("F19").Value = 868
$ws1.Range("F20").Value = 4381
$ws1.Range("F21").Value = 658
$ws1.Range("F22").Value = 198
$ws1.Range("F23").Value = 82
$ws1.Range("F26").Value = 1180
$ws1.Range("F27").Value = 89
$ws1.Range("F28").Value = 1627
$ws1.Range("F29").Value = 693
$ws1.Range("F30").Value = 882
$ws1.Range("F31").Value = 1852
$ws1.Range("F32").Value = 319
$ws1.Range("F33").Value = 2241
$ws1.Range("F34").Value = 312
$ws1.Range("F35").Value = 102
$ws1.Range("F36").Value = 1422
$ws1.Range("F39").Value = 783
$ws1.Range("F40").Value = 385
$ws1.Range("F41").Value = 4029
$ws1.Range("F44").Value = 407
$ws1.Range("F46").Value = 11
$ws1.Range("F48").Value = 161
$ws1.Range("F49").Value = 4057

# ---------------------------------------------------------------------------
# Sheet "演出" (Performance) - same kind of update on column F
# ---------------------------------------------------------------------------
$ws2.Range("F8").Value = 20
$ws2.Range("F9").Value = 20

# ---------------------------------------------------------------------------
# Sheet "本地生活" (Local life) - row 2 想去人数/最低票价 update, now sold out
# ---------------------------------------------------------------------------
$ws3.Range("F2").Value = 5120
$ws3.Range("G2").Value = "已售罄"

# ---------------------------------------------------------------------------
# Sheet "全部类型" (All types) - row 2 entirely replaced: the old "木灵动漫"
# local-life listing is gone from the combined feed and is replaced by the
# "萧敬腾·张泽" musical listing (matching 演出!B2:I2).
# ---------------------------------------------------------------------------
$ws4.Range("B2").Value = "'2024-06-28"
$ws4.Range("B2").Style = "Normal"
$ws4.Range("C2").Value = "杭州·萧敬腾·张泽领衔原创音乐剧《胭脂扣》"
$ws4.Range("D2").Value = "浙江省杭州市下城区武林广场29号 杭州剧院"
$ws4.Range("E2").Value = "2024.06.28 19:30-06.29 22:00"
$ws4.Range("F2").Value = 5
$ws4.Range("G2").Value = 480
$ws4.Range("H2").Value = "https://show.bilibili.com/platform/detail.html?id=84903"
$ws4.Range("I2").Value = "//i1.hdslb.com/bfs/openplatform/202404/NIKtJGpX1714014020771.jpeg"

# ---------------------------------------------------------------------------
# Sheet "全部类型" (All types) - remaining 想去人数 (column F) updates
# ---------------------------------------------------------------------------
$ws4.Range("F4").Value = 1500
$ws4.Range("F5").Value = 837
$ws4.Range("F6").Value = 884
$ws4.Range("F9").Value = 20
$ws4.Range("F11").Value = 5385
$ws4.Range("F12").Value = 559
$ws4.Range("F13").Value = 7450
$ws4.Range("F16").Value = 1133
$ws4.Range("F17").Value = 868
$ws4.Range("F18").Value = 4381
$ws4.Range("F19").Value = 658
$ws4.Range("F20").Value = 198
$ws4.Range("F21").Value = 82
$ws4.Range("F25").Value = 1180
$ws4.Range("F26").Value = 89
$ws4.Range("F27").Value = 1627
$ws4.Range("F28").Value = 1852
$ws4.Range("F29").Value = 319
$ws4.Range("F30").Value = 2241
$ws4.Range("F37").Value = 783
$ws4.Range("F40").Value = 385
$ws4.Range("F41").Value = 4029
$ws4.Range("F45").Value = 407
$ws4.Range("F48").Value = 161
$ws4.Range("F49").Value = 4057
